$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 165, shifting the existing rows 165-261 down
# to 166-262 (this also extends the used range to A1:R262, matching the
# dimension change in the diff, and carries the date-style formatting on
# column D down with the shifted rows).
$ws.Rows.Item(165).Insert()

# Populate the newly inserted row 165 with the new weekly record. Columns
# A,B,C,E,F,G,H,N,Q,R hold the same constant values as every other row in
# this table (market/region/category/unit/classification metadata).
$ws.Range("A165").Value = 4
$ws.Range("B165").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C165").Value = "Los Lagos"
$ws.Range("D165").Value = 44582
$ws.Range("E165").Value = 10
$ws.Range("F165").Value = 100112008
$ws.Range("G165").Value = "Coliflor"
$ws.Range("H165").Value = "Sin especificar"
$ws.Range("I165").Value = "Primera"
$ws.Range("J165").Value = 700
$ws.Range("K165").Value = 1500
$ws.Range("L165").Value = 1500
$ws.Range("M165").Value = 1500
$ws.Range("N165").Value = "`$/unidad"
$ws.Range("O165").Value = "Región Metropolitana"
$ws.Range("P165").Value = 1500
$ws.Range("Q165").Value = 1
$ws.Range("R165").Value = "Hortaliza"
